$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 13 (old rows 13-23 shift down to 15-25)
$ws.Rows("13:14").Insert()

# Clear the inherited column-A formatting/cell on the two new rows (no A cell in target)
$ws.Range("A13:A14").Clear()

# Apply column B / C body styles (style 2 / style 3) to the new rows by copying format from row 10
$ws.Range("B10").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10: Objectives text (B10/C10) changes
$ws.Range("B10").Value = 'Capacitar o aluno a interpretar e entender a natureza e a origem da estrutura e sua influência nas propriedades dos materiais cerâmicos.'
$ws.Range("C10").Value = 'Capacitar o aluno a interpretar e entender a natureza e a origem da estrutura e sua influência nas propriedades dos materiais cerâmicos.'

# New row 13: Fernando Vernilli (B/C only, A stays blank)
$ws.Range("B13").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("C13").Value = '5983729 - Fernando Vernilli Junior'

# New row 14: Sebastiao Ribeiro (B/C only, A stays blank)
$ws.Range("B14").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C14").Value = '1922320 - Sebastiao Ribeiro'

# Row 15 (was 13): Programa resumido -- B/C content changes to short syllabus text
$ws.Range("B15").Value = '1.Introdução, características dos sólidos cerâmicos, propriedades de cerâmicas'
$ws.Range("C15").Value = '1.Introdução, características dos sólidos cerâmicos, propriedades de cerâmicas'

# Row 17 (was 15): Programa -- B/C content changes to full syllabus text
$ws.Range("B17").Value = '1.Introdução1.1A indústria cerâmica1.2. Processamento de cerâmicas1.3. Produtos cerâmicos.2Características dos sólidos cerâmicos2.1Estruturas dos Cristais2.2Estruturas dos Vidros2.3Imperfeições estruturais2.4Superfície, interface e contorno de grãos2.5Mobilidade Atômica3Propriedades de Cerâmicas3.1Propriedades Térmicas3.2Propriedades Ópticas3.3Deformação Plástica, fluxo viscoso e fluência3.4Elasticidade, inelasticidade e resistência3.5Tensões térmicas e composicionais'
$ws.Range("C17").Value = '1.Introdução1.1A indústria cerâmica1.2. Processamento de cerâmicas1.3. Produtos cerâmicos.2Características dos sólidos cerâmicos2.1Estruturas dos Cristais2.2Estruturas dos Vidros2.3Imperfeições estruturais2.4Superfície, interface e contorno de grãos2.5Mobilidade Atômica3Propriedades de Cerâmicas3.1Propriedades Térmicas3.2Propriedades Ópticas3.3Deformação Plástica, fluxo viscoso e fluência3.4Elasticidade, inelasticidade e resistência3.5Tensões térmicas e composicionais'

# Row 20 (was 18): Metodo -- B/C content changes to evaluation method text
$ws.Range("B20").Value = 'Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1'
$ws.Range("C20").Value = 'Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1'

# Row 21 (was 19): Criterio -- B/C content changes to weighted average text
$ws.Range("B21").Value = 'A nota final será a média ponderada das provas escritas (80% da nota final) e das listas de exercícios e relatórios (20% da nota final).'
$ws.Range("C21").Value = 'A nota final será a média ponderada das provas escritas (80% da nota final) e das listas de exercícios e relatórios (20% da nota final).'

# Row 22 (was 20): Norma de recuperacao -- B/C content changes to recovery exam text
$ws.Range("B22").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). NF = (MP + PR)/2. NF igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado.'
$ws.Range("C22").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). NF = (MP + PR)/2. NF igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado.'

# Row 23 (was 21): Bibliografia -- B/C content changes to bibliography text
$ws.Range("B23").Value = '1. KINGERY, W. D.; BOWEN, H. K.; UHLMANN, D. R. Introduction of ceramics New York: John Wiley, c1976
2. R.W. CAHN; P.HAASEN; E.J. KRAMER. Materials Science and Technology: A      Comprehensive Treatment. Weinheim: Wiley-VCH, c2005BERGERON, CLIFTON G.; RISBUD, SUBHASH H. Introduction to phase equilibria in ceramics. Westerville: The American Ceramic Society, 1984.BROOK, R. J. Processing of ceramics. R. W. Cahn; P. Haasen; E. J. Kramer. Weinheim: VCH, 1996.LEVIN, ERNEST M. Phase diagrams for ceramics. Ohio: The American Ceramic, 1964.R.C. BRADT; D.P.H.HASSELMAN; D. MUNZ; M.SAKAI; V.YASHEVCHENKO  Fracture mechanics of ceramics: r-curve behavior, toughness determination, and thermal shock.. New York: Plenum, 1996.R.C. BRADT; D.P.H.HASSELMAN; D. MUNZ; M.SAKAI; V.YASHEVCHENKO  Fracture mechanics of ceramics: fatigue, composites, and high-temperature behaviour.. New York: Plenum, 1996.REED, JAMES S. Principles of ceramics processing. New YorK: John Wiley, 1995.TOMPSON, D.P., ED. Engineering ceramics: fabrication science & technology.  London: The Institute of Materials, 1993.R.C. BRADT; D.P.H.HASSELMAN; D. MUNZ; M.SAKAI; V.YASHEVCHENKO  Fracture mechanics of ceramics: fatigue, composites, and high-temperature behaviour.. New York: Plenum, 1996.CHIANG, YET-MING; BIRNIE III, DUNBAR P.; KINGERY, W.DAVID. Physical ceramics: principles for ceramic science and engineering. New York: John Wiley, 1997.MENCIK, JAROSLAV. Strength and fracture of glass and ceramics.  Amsterdam: Elsevier, 1992.'
$ws.Range("C23").Value = '1. KINGERY, W. D.; BOWEN, H. K.; UHLMANN, D. R. Introduction of ceramics New York: John Wiley, c1976
2. R.W. CAHN; P.HAASEN; E.J. KRAMER. Materials Science and Technology: A      Comprehensive Treatment. Weinheim: Wiley-VCH, c2005BERGERON, CLIFTON G.; RISBUD, SUBHASH H. Introduction to phase equilibria in ceramics. Westerville: The American Ceramic Society, 1984.BROOK, R. J. Processing of ceramics. R. W. Cahn; P. Haasen; E. J. Kramer. Weinheim: VCH, 1996.LEVIN, ERNEST M. Phase diagrams for ceramics. Ohio: The American Ceramic, 1964.R.C. BRADT; D.P.H.HASSELMAN; D. MUNZ; M.SAKAI; V.YASHEVCHENKO  Fracture mechanics of ceramics: r-curve behavior, toughness determination, and thermal shock.. New York: Plenum, 1996.R.C. BRADT; D.P.H.HASSELMAN; D. MUNZ; M.SAKAI; V.YASHEVCHENKO  Fracture mechanics of ceramics: fatigue, composites, and high-temperature behaviour.. New York: Plenum, 1996.REED, JAMES S. Principles of ceramics processing. New YorK: John Wiley, 1995.TOMPSON, D.P., ED. Engineering ceramics: fabrication science & technology.  London: The Institute of Materials, 1993.R.C. BRADT; D.P.H.HASSELMAN; D. MUNZ; M.SAKAI; V.YASHEVCHENKO  Fracture mechanics of ceramics: fatigue, composites, and high-temperature behaviour.. New York: Plenum, 1996.CHIANG, YET-MING; BIRNIE III, DUNBAR P.; KINGERY, W.DAVID. Physical ceramics: principles for ceramic science and engineering. New York: John Wiley, 1997.MENCIK, JAROSLAV. Strength and fracture of glass and ceramics.  Amsterdam: Elsevier, 1992.'

# Row 24 (new): Requisitos label
$ws.Range("A24").Value = 'Requisitos:'

# Row 25 (new): LOM3013 requirement text B/C
$ws.Range("B25").Value = 'LOM3013 -  Ciência dos Materiais  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOM3013 -  Ciência dos Materiais  (Requisito fraco)
'
